# Apply "new TPM" data update to the Vip-Sctr LR-pairs sheet.
# The sheet originally held 3 data rows (target clusters ECs, MuSCs,
# Resolving-Mac). After the update only the Resolving-Mac row survives,
# now carrying freshly recomputed values, so we:
#   1. drop the two rows that are no longer needed (old row 2 "ECs" data
#      and old row 3 "MuSCs" data), keeping the old row 4 position's data
#      but overwriting it with the new numbers in row 2
#   2. rewrite the remaining data row's values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 and 4 entirely, leaving header (row 1) + single data row (row 2)
$ws.Rows("3:4").Delete()

# Update the Target cluster label for the remaining row
$ws.Range("D2").Value = "Resolving-Mac"

# Update the recomputed TPM-based metrics for the remaining row
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3882076666666667
$ws.Range("H2").Value = 1.164623
$ws.Range("M2").Value = 0.02199266666666666
$ws.Range("N2").Value = 0.06597799999999999
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.008537721810444444
$ws.Range("R2").Value = 0.076839496294
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
